# SCD0022-004 - Update TC_ID values and sheet name
# (Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from SCD0317 to SCD0022
$ws.Name = "SCD0022"

# Update the TC_ID column (B) values for rows 2-5 from "DGS-332" to "SCD0022-004"
$ws.Range("B2").Value2 = "SCD0022-004"
$ws.Range("B3").Value2 = "SCD0022-004"
$ws.Range("B4").Value2 = "SCD0022-004"
$ws.Range("B5").Value2 = "SCD0022-004"

# Column B needs to widen to fit the new, longer TC_ID text
$ws.Columns.Item(2).AutoFit()

# Move the selection to B6, matching where the editor left off
$ws.Range("B6").Select()
